$d = $word.ActiveDocument

# 1. Replace the hard-coded company name with the new $CURRENTUSERCOMPANY$
#    replacement variable (commit: "feat: New replacement variables OFFICE, COMPANY").
$old = "Super Duper Inc."
$new = "`$CURRENTUSERCOMPANY`$"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# 2. Remove the stray "_GoBack" bookmark that Word leaves behind after the
#    last edit position; newer Word versions no longer emit it on save.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
